# Apply cryptos list price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.334.86"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.668.17"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.53"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5311"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2648"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06369"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.95"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07843"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.533"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.670.92"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "1.896.80"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5613"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "0.0₅8146"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "26.339.69"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.721"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "197.92"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.054"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.010"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.74"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1217"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.254"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05891"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.286"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.552"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.327"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.603"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9618"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.435"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5814"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.963"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").Value = "1.074.65"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8581"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.009"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.79"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("D45").Value = "1.807.61"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.52"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.014"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4410"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.074"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05148"
$ws.Range("E51").Value = "  -0.17%  "
